# "complete the house data by coredata"
# Adds a new "sortOrder" column (H) to the houseConfig sheet with an
# int type-row and per-item sort values, then leaves the workbook's
# active sheet/selection the way the author left it (houseConfig
# active on cell F9, functionaryConfig selection parked at D4).

$wb = $excel.ActiveWorkbook

$wsHouse = $wb.Worksheets.Item("houseConfig")
$wsFunctionary = $wb.Worksheets.Item("functionaryConfig")

# New column header + type row
$wsHouse.Range("H1").Value = "sortOrder"
$wsHouse.Range("H2").Value = "int"

# Per-row sortOrder values
$wsHouse.Range("H3").Value = 1
$wsHouse.Range("H4").Value = 2
$wsHouse.Range("H5").Value = 3

# Restore the view/selection state recorded in the saved workbook:
# houseConfig becomes the active/selected sheet again, with the
# cursor left on F9; functionaryConfig's lingering selection moves
# to D4 and it is no longer the tab in focus.
$wsFunctionary.Range("D4").Select() | Out-Null
$wsHouse.Activate() | Out-Null
$wsHouse.Range("F9").Select() | Out-Null
